$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timestamp values (Excel serial date numbers) for column D,
# shifted down one "slot" as a fresh update was appended.
$valGroup1 = 44242.54844760572   # rows 2-15  (was 44242.52724028876)
$valGroup2 = 44242.52724028935   # rows 16-29 (was 44242.50600597222)
$valGroup3 = 44242.50600597222   # rows 30-43 (was 44242.48478574074)

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = $valGroup1
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = $valGroup2
}

for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = $valGroup3
}
